$d = $word.ActiveDocument

# Locate the last table in the "Invalidate Break-in" use case (the one whose
# last cell contains "Break-in has been handled") and insert the new
# "Test Plan" paragraphs right after it, before the paragraph that already
# follows the table.
$tbl = $d.Tables.Item($d.Tables.Count)
$insertionPoint = $d.Range($tbl.Range.End, $tbl.Range.End)

$CR = [char]13

# Insert in reverse order: each InsertBefore call places its text
# immediately before $insertionPoint, which stays anchored right before the
# paragraph that originally followed the table, so inserting back-to-front
# yields the correct final top-to-bottom order.
$insertionPoint.InsertBefore($CR)
$insertionPoint.InsertBefore("2. Test if password is not valid, break in has been handled" + $CR)
$insertionPoint.InsertBefore("1. Test ideal path, ensure post conditions have been met" + $CR)
$insertionPoint.InsertBefore("Test Plan:" + $CR)
$insertionPoint.InsertBefore($CR)
